# The document has two "logo" pictures that each appear twice (once in the
# default header/footer, once in the first-page header/footer):
#   - the Pearson logo, an inline picture embedded as image1.png, shown in
#     both footers with docPr/cNvPr name="image1.png"
#   - the BTEC logo, an inline picture embedded as image2.jpg, shown in
#     both headers with docPr/cNvPr name="image2.jpg"
#
# This edit renames the *display name* of each inline picture (the
# wp:docPr/@name that Word exposes as Shape.Name) so the two pairs swap
# numbers: the Pearson logo pictures become "image2.png" and the BTEC logo
# pictures become "image1.jpg".
#
# InlineShape has no settable Name property in the Word object model, so
# each picture is temporarily converted to a floating Shape (which does
# expose .Name), renamed, and converted back to an inline picture.

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

function Rename-LogoPicture($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Footers: Pearson logo, image1.png -> image2.png
Rename-LogoPicture $section.Footers.Item(1).Range "image2.png"
Rename-LogoPicture $section.Footers.Item(2).Range "image2.png"

# Headers: BTEC logo, image2.jpg -> image1.jpg
Rename-LogoPicture $section.Headers.Item(1).Range "image1.jpg"
Rename-LogoPicture $section.Headers.Item(2).Range "image1.jpg"
